$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: phase / plan-type header labels for the new summary block ---
$ws.Range("E17").Value = "DHMO"
$ws.Range("F17").Value = "AHMO"
$ws.Range("G17").Value = "DPPO"
$ws.Range("H17").Value = "APPO"

# --- Row 19: ape2e ---
$ws.Range("A19").Value = "ape2e"
$ws.Range("C19").Value = 72
$ws.Range("H19").Formula = "=C19"

# --- Row 20: bkp ---
$ws.Range("A20").Value = "bkp"
$ws.Range("C20").Value = 136
$ws.Range("G20").Formula = "=C20"

# --- Row 21: content ---
$ws.Range("A21").Value = "content"
$ws.Range("C21").Value = 63
$ws.Range("E21").Formula = "=C21"
$ws.Range("G21").Formula = "=C21"
$ws.Range("H21").Formula = "=C21"

# --- Row 22: deps ---
$ws.Range("A22").Value = "deps"
$ws.Range("C22").Value = 38
$ws.Range("E22").Formula = "=C22"
$ws.Range("F22").Formula = "=C22"
$ws.Range("G22").Formula = "=C22"
$ws.Range("H22").Formula = "=C22"

# --- Row 23: dpe2e ---
$ws.Range("A23").Value = "dpe2e"
$ws.Range("C23").Value = 108
$ws.Range("G23").Formula = "=C23"

# --- Row 24: e2eAP ---
$ws.Range("A24").Value = "e2eAP"
$ws.Range("C24").Value = 110
$ws.Range("E24").Formula = "=C24"
$ws.Range("F24").Formula = "=C24"
$ws.Range("G24").Formula = "=C24"
$ws.Range("H24").Formula = "=C24"

# --- Row 25: facs ---
$ws.Range("A25").Value = "facs"
$ws.Range("C25").Value = 20
$ws.Range("E25").Formula = "=C25"
$ws.Range("F25").Formula = "=C25"

# --- Row 26: pays ---
$ws.Range("A26").Value = "pays"
$ws.Range("C26").Value = 17
$ws.Range("E26").Formula = "=C26"
$ws.Range("F26").Formula = "=C26"
$ws.Range("G26").Formula = "=C26"
$ws.Range("H26").Formula = "=C26"

# --- Row 27: pdf ---
$ws.Range("A27").Value = "pdf"
$ws.Range("C27").Value = 8
$ws.Range("G27").Formula = "=C27"
$ws.Range("H27").Formula = "=C27"

# --- Row 28: pers ---
$ws.Range("A28").Value = "pers"
$ws.Range("C28").Value = 122
$ws.Range("E28").Formula = "=C28"
$ws.Range("F28").Formula = "=C28"
$ws.Range("G28").Formula = "=C28"
$ws.Range("H28").Formula = "=C28"

# --- Row 29: shop ---
$ws.Range("A29").Value = "shop"
$ws.Range("C29").Value = 41
$ws.Range("E29").Formula = "=C29"
$ws.Range("F29").Formula = "=C29"
$ws.Range("G29").Formula = "=C29"
$ws.Range("H29").Formula = "=C29"

# --- Row 31: totals across the new block ---
$ws.Range("C31").Formula = "=SUM(C19:C29)"
$ws.Range("E31").Formula = "=SUM(E19:E29)"
$ws.Range("F31").Formula = "=SUM(F19:F29)"
$ws.Range("G31").Formula = "=SUM(G19:G29)"
$ws.Range("H31").Formula = "=SUM(H19:H29)"
$ws.Range("J31").Formula = "=SUM(E31:H31)"

# --- Row 32: secondary snapshot totals ---
$ws.Range("F32").Value = 348
$ws.Range("H32").Value = 471
$ws.Range("J32").Formula = "=SUM(E32:H32)"

# --- View state: scroll/selection as left by the author ---
$ws.Range("C22").Select()
